$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 3280.3076
$ws.Range("J62").Value = 8500
$ws.Range("L62").Value = 8500
$ws.Range("N62").Value = -9748

$ws.Range("H65").Value = 3280.3076
$ws.Range("J65").Value = 8500
$ws.Range("L65").Value = 42500
$ws.Range("N65").Value = -48740

$ws.Range("H70").Value = 2560.7778
$ws.Range("I70").Value = 1991.1578
$ws.Range("J70").Value = 3913.625
$ws.Range("K70").Value = 5973.4734
$ws.Range("L70").Value = 11740.875
$ws.Range("M70").Value = -5703.4734
$ws.Range("N70").Value = -12280.875

$ws.Range("H73").Value = 2560.7778
$ws.Range("I73").Value = 1991.1578
$ws.Range("J73").Value = 3913.625
$ws.Range("K73").Value = 5973.4734
$ws.Range("L73").Value = 11740.875
$ws.Range("M73").Value = -5037.4734
$ws.Range("N73").Value = -13612.875

$ws.Range("H107").Value = 2360.3333
$ws.Range("I107").Value = 2903.5715
$ws.Range("J107").Value = 1599.8
$ws.Range("K107").Value = 2903.5715
$ws.Range("L107").Value = 1599.8
$ws.Range("M107").Value = -983.5715
$ws.Range("N107").Value = -5439.8

$ws.Range("H113").Value = 6501.222
$ws.Range("I113").Value = 3215
$ws.Range("J113").Value = 18003
$ws.Range("K113").Value = 3215
$ws.Range("L113").Value = 18003
$ws.Range("M113").Value = 39
$ws.Range("N113").Value = -24511

$ws.Range("H116").Value = 632893.5
$ws.Range("I116").Value = 1253548.8
$ws.Range("J116").Value = 12238.25
$ws.Range("K116").Value = 1253548.8
$ws.Range("L116").Value = 12238.25
$ws.Range("M116").Value = -1250106.8
$ws.Range("N116").Value = -19122.25

$ws.Range("H141").Value = 96936.91
$ws.Range("I141").Value = 106625
$ws.Range("K141").Value = 319875
$ws.Range("M141").Value = -314695

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 3416.5
$ws.Range("I122").Value = 1761.7142
$ws.Range("J122").Value = 15000
$ws.Range("K122").Value = 5285.142599999999
$ws.Range("L122").Value = 45000
$ws.Range("M122").Value = -2835.142599999999
$ws.Range("N122").Value = -49900

$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()

$ws.Range("H133").Value = 20400.285
$ws.Range("J133").Value = 20400.285
$ws.Range("L133").Value = 20400.285
$ws.Range("N133").Value = -25460.285

$ws.Range("H137").Value = 40980
$ws.Range("J137").Value = 40980
$ws.Range("L137").Value = 40980
$ws.Range("N137").Value = -51180

$ws.Range("H139").Value = 40823.516
$ws.Range("J139").Value = 40823.516
$ws.Range("L139").Value = 40823.516
$ws.Range("N139").Value = -51103.516

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H132").Value = 50845
$ws.Range("J132").Value = 50845
$ws.Range("L132").Value = 50845
$ws.Range("N132").Value = -60965

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("I31").Value = 1109.6154
$ws.Range("J31").Value = 5842.75
$ws.Range("K31").Value = 1109.6154
$ws.Range("L31").Value = 5842.75
$ws.Range("M31").Value = -814.6153999999999
$ws.Range("N31").Value = -6432.75

$ws.Range("I34").Value = 1109.6154
$ws.Range("J34").Value = 5842.75
$ws.Range("K34").Value = 1109.6154
$ws.Range("L34").Value = 5842.75
$ws.Range("M34").Value = -907.6153999999999
$ws.Range("N34").Value = -6246.75

$ws.Range("H122").Value = 4953
$ws.Range("I122").Value = 1604
$ws.Range("J122").Value = 15000
$ws.Range("K122").Value = 4812
$ws.Range("L122").Value = 45000
$ws.Range("M122").Value = -2362
$ws.Range("N122").Value = -49900

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 88.63158
$ws.Range("I12").Value = 28.6
$ws.Range("J12").Value = 110.07143
$ws.Range("K12").Value = 85.80000000000001
$ws.Range("L12").Value = 330.21429
$ws.Range("M12").Value = 87.19999999999999
$ws.Range("N12").Value = -676.21429

$ws.Range("H87").Value = 3230.4285
$ws.Range("I87").Value = 871.3333
$ws.Range("K87").Value = 2613.9999
$ws.Range("M87").Value = -1365.9999

$ws.Range("H90").Value = 3230.4285
$ws.Range("I90").Value = 871.3333
$ws.Range("K90").Value = 7841.9997
$ws.Range("M90").Value = -1601.9997

$ws.Range("H98").Value = 475
$ws.Range("I98").Value = 475
$ws.Range("K98").Value = 1425
$ws.Range("M98").Value = 73

$ws.Range("H131").Value = 710.6
$ws.Range("J131").Value = 798.51807
$ws.Range("L131").Value = 2395.55421
$ws.Range("N131").Value = -12475.55421

$ws.Range("H133").Value = 6540
$ws.Range("I133").Value = 6975
$ws.Range("K133").Value = 20925
$ws.Range("M133").Value = -15865

$ws.Range("H134").Value = 6763.846
$ws.Range("I134").Value = 10957.5
$ws.Range("J134").Value = 4900
$ws.Range("K134").Value = 32872.5
$ws.Range("L134").Value = 14700
$ws.Range("M134").Value = -27802.5
$ws.Range("N134").Value = -24840

$ws.Range("H137").Value = 1393.3334
$ws.Range("I137").Value = 1393.3334
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 4180.0002
$ws.Range("L137").Value = 0
$ws.Range("M137").Value = 919.9997999999996
$ws.Range("N137").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1912.5
$ws.Range("I102").Value = 1250.12
$ws.Range("J102").Value = 5224.4
$ws.Range("K102").Value = 1250.12
$ws.Range("L102").Value = 5224.4
$ws.Range("M102").Value = 371.8800000000001
$ws.Range("N102").Value = -8468.4

$ws.Range("H113").Value = 1706.3529
$ws.Range("I113").Value = 1708.3077
$ws.Range("J113").Value = 1700
$ws.Range("K113").Value = 1708.3077
$ws.Range("L113").Value = 1700
$ws.Range("M113").Value = 461.6922999999999
$ws.Range("N113").Value = -6040

$ws.Range("H122").Value = 15000
$ws.Range("I122").Value = 5000
$ws.Range("J122").Value = 18333.334
$ws.Range("K122").Value = 15000
$ws.Range("L122").Value = 55000.00199999999
$ws.Range("M122").Value = -12550
$ws.Range("N122").Value = -59900.00199999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 5672.727
$ws.Range("I122").Value = 1566.6666
$ws.Range("J122").Value = 10600
$ws.Range("K122").Value = 4699.9998
$ws.Range("L122").Value = 31800
$ws.Range("M122").Value = -2249.9998
$ws.Range("N122").Value = -36700

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 10659.25
$ws.Range("I45").Value = 8784.5
$ws.Range("J45").Value = 12534
$ws.Range("K45").Value = 8784.5
$ws.Range("L45").Value = 12534
$ws.Range("M45").Value = -8293.5
$ws.Range("N45").Value = -13516

$ws.Range("H80").Value = 39466
$ws.Range("J80").Value = 39466
$ws.Range("L80").Value = 39466
$ws.Range("N80").Value = -41462

$ws.Range("H83").Value = 39466
$ws.Range("J83").Value = 39466
$ws.Range("L83").Value = 118398
$ws.Range("N83").Value = -128382

$ws.Range("H122").Value = 3585.7896
$ws.Range("I122").Value = 1414.0714
$ws.Range("J122").Value = 9666.6
$ws.Range("K122").Value = 4242.2142
$ws.Range("L122").Value = 28999.8
$ws.Range("M122").Value = -1792.2142
$ws.Range("N122").Value = -33899.8
